$d = $word.ActiveDocument

# 1. "Cancelar y agregar uno o varios ítems del pedido." -> "Cancelar uno o varios ítems del pedido."
$d.Content.Find.Execute(
    "Cancelar y agregar uno o varios ítems del pedido.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cancelar uno o varios ítems del pedido.",
    2)

# 2. "Que se cancele y agreguen nuevos ítems si el cliente lo desea." ->
#    "Que se cancele ítems del pedido si el cliente lo desea."
$d.Content.Find.Execute(
    "Que se cancele y agreguen nuevos ítems si el cliente lo desea.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Que se cancele ítems del pedido si el cliente lo desea.",
    2)

# 3. Update the "Objetivo Funcionalidad" sentence describing the cancel flow.
#    Done in two pieces so the _GoBack bookmark (which sits between the old
#    second and third runs) is not spanned/removed by either replacement;
#    afterwards it naturally ends up trailing at the end of the paragraph.
$d.Content.Find.Execute(
    "Antes de generar la factura final si el cliente desea cancelar uno de los ítems del pedido y/o agregar uno nuevo, con el botón cancelar en factura se regresa y se realiza",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Antes de generar la factura final si el cliente desea cancelar uno de los ítems del pedido con el botón cancelar en factura se regresa y se realiza la operación de cancelar un ítem del pedido haciendo clic sobra el botón con logo de X.",
    2)

$d.Content.Find.Execute(
    " la operación que desea el cliente.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2)
